$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.559.28'
$ws.Range('D3').Value = '1.741.29'
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '246.11'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +5.20%  '
$ws.Range('E6').Value = '  -0.09%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4803'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +4.53%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2685'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +4.55%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06252'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +2.31%  '
$ws.Range('D10').Value = '1.742.31'
$ws.Range('E10').Value = '  +4.55%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07142'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +2.72%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '15.78'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +8.24%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.6209'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +9.67%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.526'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +4.56%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '77.31'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +3.49%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.001'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -0.08%  '
$ws.Range('D17').Value = '26.564.35'
$ws.Range('E18').Value = '  -0.15%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000006900'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +3.53%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.75'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +3.66%  '
$ws.Range('D21').Value = '1.966.45'
$ws.Range('E21').Value = '  +4.52%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.608'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +4.40%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.869'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +2.00%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.362'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +3.15%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '135.99'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.09%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '15.36'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +3.13%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.808'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +6.63%  '
$ws.Range('E28').Value = '  +3.23%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '106.87'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +2.74%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.004'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +1.27%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.740'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +4.10%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.07874'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +1.62%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04571'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +6.95%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.619'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -0.34%  '
$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9999'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +5.88%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.6364'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +6.74%  '
$ws.Range('E37').Value = '  +1.29%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '113.18'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +11.84%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.431'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -3.42%  '
$ws.Range('E40').Value = '  +8.99%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.005'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +0.29%  '
$ws.Range('B42').Value = 'VeChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.01512'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +3.62%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.731'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +16.80%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.3907'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +6.08%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '6.769'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +10.96%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.1211'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +9.99%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.05331'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +1.43%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.889'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +7.34%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '30.75'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +3.98%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.258'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +5.47%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.3442'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +5.00%  '
